$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents (values + formatting) of column B and column C ---
# Before the edit: B held the day-index (Timp(zile)) and C held the calendar
# date (Data calendaristica). After the edit those two columns trade places:
# B now holds the date and C now holds the day-index. Using Copy(Destination)
# carries formatting (fill/number-format/alignment) along with the value, so
# a 3-way swap through a scratch column reproduces both the content and the
# style-index swap seen in the diff (B1/C1 and B2:B196/C2:C196).
$scratchCol = "BZ"
$lastRow = 196

$srcB = $ws.Range("B1:B$lastRow")
$srcC = $ws.Range("C1:C$lastRow")
$scratch = $ws.Range("${scratchCol}1:${scratchCol}$lastRow")

$srcB.Copy($scratch)
$srcC.Copy($srcB)
$scratch.Copy($srcC)
$scratch.Clear()

# --- Column widths: column B (now holding dates) needs to be widened; the
# rest of the bestFit block (C:Q) keeps its previous width. ---
$ws.Columns("B").ColumnWidth = 16.3

# --- Selection / view state: the whole column B is now selected (instead of
# cell A5), and the sheet view no longer pins a scrolled-down top-left cell. ---
$ws.Range("A1").Select() | Out-Null
$ws.Columns("B").Select() | Out-Null
